$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F2, F4, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1120
$ws1.Range("F4").Value = 1816
$ws1.Range("F6").Value = 436

# Sheet "全部类型" (sheet4): update F2, F4, F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1120
$ws4.Range("F4").Value = 1816
$ws4.Range("F7").Value = 436
